$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN  ['Bernese Mountain Dog']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Files-tab Cypher query in B4: the corrected script drops the
# redundant `File Type` and `Breed` columns from the RETURN clause.
$ws.Range("B4").Value = $newFilesQuery

# The row was sized to fit the old (longer) wrapped text; with two fewer
# lines the row now needs less height.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moves to the edited cell.
[void]$ws.Range("B4").Select()
